# Commit: "#5: insurance, claim, debt, investment done"
#
# Adds property_category / category / date / legislator_name / legislator_id /
# source_file / index columns (plus a couple of reshuffled header/label cells)
# to the 保險 (insurance), 債權 (claim) and 事業投資 (investment) sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "具有相當價值之財產" (property of considerable value) -- the
# property_category label for the jewellery/antiques rows was "otherbonds";
# it is recategorised to "antique". (F2:F18, all rows of the sheet.)
# ---------------------------------------------------------------------------
$ws0 = $wb.Worksheets.Item("具有相當價值之財產")
for ($r = 2; $r -le 18; $r++) {
    if ($ws0.Cells.Item($r, 6).Value2 -eq "otherbonds") {
        $ws0.Cells.Item($r, 6).Value = "antique"
    }
}

# ---------------------------------------------------------------------------
# Sheet "保險" (insurance) -- was A1:E8, becomes A1:K8
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("保險")

# Header row: B/C/D/E had been storing literal row-2 data instead of field
# names -- fix them up, then append the new trailing field-name columns.
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

$rows6 = @(
  @{ R=2; Idx=312; Company="國泰人壽"; Name="利率變動型年金保險（甲型）"; Owner="陳錦錠" },
  @{ R=3; Idx=313; Company="國泰人壽"; Name="金歡喜110養老保險";       Owner="陳錦錠" },
  @{ R=4; Idx=314; Company="國泰人壽"; Name="有GO讃養老保險";         Owner="陳錦錠" },
  @{ R=5; Idx=315; Company="國泰人壽"; Name="金好鑽養老保險";         Owner="陳錦錠" },
  @{ R=6; Idx=317; Company="國泰人壽"; Name="達康101終身保險";        Owner="陳錦錠" },
  @{ R=7; Idx=318; Company="新光人壽"; Name="長樂終身險";             Owner="張慶忠" },
  @{ R=8; Idx=319; Company="國泰人壽"; Name="金美利美元養老保險";     Owner="陳錦錠" }
)

foreach ($row in $rows6) {
    $r = $row.R
    $ws.Cells.Item($r, 2).Value = $row.Company          # B: company
    $ws.Cells.Item($r, 3).Value = $row.Name              # C: name
    $ws.Cells.Item($r, 4).Value = $row.Owner             # D: owner
    $ws.Cells.Item($r, 5).Value = "insurance"            # E: property_category
    $ws.Cells.Item($r, 6).Value = "normal"               # F: category
    $ws.Cells.Item($r, 7).Value = "'2013-12-12"          # G: date (leading ' keeps it text, not an auto-converted date)
    $ws.Cells.Item($r, 8).Value = "張慶忠"                # H: legislator_name
    $ws.Cells.Item($r, 9).Value = 1347                   # I: legislator_id
    $ws.Cells.Item($r, 10).Value = "tmpe4561"            # J: source_file
    $ws.Cells.Item($r, 11).Value = $row.Idx              # K: index
}

# ---------------------------------------------------------------------------
# Sheet "債權" (claim) -- was A1:G7, becomes A1:N7
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("債權")

$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "debtor"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

$rows7 = @(
  @{ R=2; Idx=324; Total=180350000; Debtor="漢寶開發建設(股)公司新北市中和區安樂路217巷20巷8號"; RegDate="98年10月23日";  Reason="借款" },
  @{ R=3; Idx=325; Total=44665413;  Debtor="基鴻建設股份有限公司新北市中和區建一路148號2樓.";        RegDate="100年08月15日"; Reason="股東往來" },
  @{ R=4; Idx=326; Total=140700000; Debtor="漢龍營造股份有限公司新北市中和區安樂路217巷20弄8號";     RegDate="98年11月23日"; Reason="借款" },
  @{ R=5; Idx=327; Total=177126000; Debtor="漢寶開發建設(股)公司新北市中和區安樂路217巷20巷8號";     RegDate="102年12月12日"; Reason="代收土地款" },
  @{ R=6; Idx=328; Total=66700000;  Debtor="漢堡開發建設(股)公司新北市中和區建一路148號8樓";        RegDate="100年02月22日"; Reason="借款" },
  @{ R=7; Idx=329; Total=180588000; Debtor="資信建設股份有限公司新北市中和區安樂路217巷20弄8號";     RegDate="98年12月06日"; Reason="借款" }
)

foreach ($row in $rows7) {
    $r = $row.R
    $ws.Cells.Item($r, 2).Value = "未兌現支票"           # B: species
    $ws.Cells.Item($r, 3).Value = "張慶忠"                # C: owner (row3 overridden below)
    $ws.Cells.Item($r, 4).Value = $row.Debtor             # D: debtor
    $ws.Cells.Item($r, 5).Value = $row.Total              # E: total
    $ws.Cells.Item($r, 6).Value = $row.RegDate            # F: register_date
    $ws.Cells.Item($r, 7).Value = $row.Reason             # G: register_reason
    $ws.Cells.Item($r, 8).Value = "claim"                 # H: property_category
    $ws.Cells.Item($r, 9).Value = "normal"                # I: category
    $ws.Cells.Item($r, 10).Value = "'2013-12-12"          # J: date (leading ' keeps it text, not an auto-converted date)
    $ws.Cells.Item($r, 11).Value = "張慶忠"               # K: legislator_name
    $ws.Cells.Item($r, 12).Value = 1347                   # L: legislator_id
    $ws.Cells.Item($r, 13).Value = "tmpe4561"             # M: source_file
    $ws.Cells.Item($r, 14).Value = $row.Idx               # N: index
}
# Row 3 (index 325) owner is 陳錦錠, not 張慶忠 -- fix it up.
$ws.Cells.Item(3, 3).Value = "陳錦錠"

# ---------------------------------------------------------------------------
# Sheet "事業投資" (investment) -- was A1:G3, becomes A1:N3
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("事業投資")

$ws.Range("B1").Value = "owner"
$ws.Range("C1").Value = "company"
$ws.Range("D1").Value = "address"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

$rows8 = @(
  @{ R=2; Idx=339; Owner="張慶忠";  Company="漢禧建設有限公司";       Address="新北市安樂路217巷20弄8號"; Total=17500000; RegDate="86年03月18日"; Reason="股金" },
  @{ R=3; Idx=340; Owner="張慶忠.";  Company="弘基消防安全設備有限公司"; Address="新北市景平路7821號5樓";    Total=2000000;  RegDate="92年03月28日"; Reason="股金" }
)

foreach ($row in $rows8) {
    $r = $row.R
    $ws.Cells.Item($r, 2).Value = $row.Owner              # B: owner
    $ws.Cells.Item($r, 3).Value = $row.Company            # C: company
    $ws.Cells.Item($r, 4).Value = $row.Address            # D: address
    $ws.Cells.Item($r, 5).Value = $row.Total              # E: total
    $ws.Cells.Item($r, 6).Value = $row.RegDate            # F: register_date
    $ws.Cells.Item($r, 7).Value = $row.Reason             # G: register_reason
    $ws.Cells.Item($r, 8).Value = "investment"            # H: property_category
    $ws.Cells.Item($r, 9).Value = "normal"                # I: category
    $ws.Cells.Item($r, 10).Value = "'2013-12-12"          # J: date (leading ' keeps it text, not an auto-converted date)
    $ws.Cells.Item($r, 11).Value = "張慶忠"               # K: legislator_name
    $ws.Cells.Item($r, 12).Value = 1347                   # L: legislator_id
    $ws.Cells.Item($r, 13).Value = "tmpe4561"             # M: source_file
    $ws.Cells.Item($r, 14).Value = $row.Idx               # N: index
}
